$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh - GitHub Actions scheduled update
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.607.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.84%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.248.64'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '184.05'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.86%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("E8").Value = '  +1.36%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.246.27'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.84%  '
$ws.Range("E10").Value = '  -3.16%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.58'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -2.37%  '
$ws.Range("E12").Value = '  -1.22%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.799.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.16%  '
$ws.Range("E14").Value = '  +0.05%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.66'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -3.12%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.637.10'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.78%  '
$ws.Range("E17").Value = '  -1.74%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.213.36'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '395.30'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.57'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.28%  '
$ws.Range("E23").Value = '  +0.00%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.516'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("E26").Value = '  -2.73%  '
$ws.Range("E27").Value = '  -1.66%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.64%  '
$ws.Range("E29").Value = '  -0.38%  '
$ws.Range("E30").Value = '  -2.23%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.55'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.00%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.66'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.31%  '
$ws.Range("E33").Value = '  -2.02%  '
$ws.Range("E34").Value = '  -2.17%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '161.35'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("E37").Value = '  -4.07%  '
$ws.Range("E38").Value = '  +1.41%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '26.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.25%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.807'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.86%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '4.54'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.43%  '
$ws.Range("E42").Value = '  -4.47%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.47'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.53%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.0687'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '40.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.73%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.615.06'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.95%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.70'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -3.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '334.94'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0278'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.24%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.97%  '
$ws.Range("E51").Value = '  -0.46%  '
